# Edit: "Black Holes" -> "Chemistry" content rewrite
# Title, author name/email, intro body, summary body all rewritten with new
# chemistry-themed copy; intro body gets substantially expanded with new
# sentences/paragraph breaks; a trailing empty paragraph is appended.

$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Title
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unraveling the Enigma of Black Holes", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "Unveiling the Magic of Chemistry: Exploring the World of Atoms, Molecules, and Chemical Reactions",
    2) | Out-Null

# -----------------------------------------------------------------------
# 2) Author name
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    " Sophia Carter", $false, $false, $false, $false,
    $false, $true, 1, $false,
    " Alex Spencer",
    2) | Out-Null

# -----------------------------------------------------------------------
# 3) Author email (two runs: local/domain part, and tld part)
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "sophiacarter@astro", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "alexspencer@gmail",
    2) | Out-Null

$d.Content.Find.Execute(
    "edu", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "com",
    2) | Out-Null

# -----------------------------------------------------------------------
# 4) Intro body paragraph - full replace (sentence count changes, so the
#    run layout changes too: some runs are dropped, many new ones added).
# -----------------------------------------------------------------------
$introXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00653B50" w:rsidRDefault="005C06F8">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>Step into the fascinating realm of chemistry, where we delve into the microscopic world of atoms, molecules, and intriguing chemical reactions</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Chemistry serves as the foundation for understanding the composition, properties, and changes that occur in matter around us</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> It plays a crucial role in various fields, from medicine to materials science, unraveling the complexities of how substances interact with each other</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>In this exciting journey into the world of chemistry, we'll embark on the analysis of chemical reactions, discovering how elements combine and rearrange to form new substances</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Unveil the intricacies of periodic trends and properties, leading us to appreciate the structure and behavior of molecules</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> We'll explore the magic behind everyday phenomena, shedding light on the chemistry of cooking and the intricate mechanisms responsible for cellular respiration</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>As we delve deeper into the world of atoms and molecules, we'll encounter the awe-inspiring beauty of chemistry</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Chemistry is a harmonious symphony of colors, reactions, and energy transformations, waiting to be appreciated by curious minds</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> With each new concept unravelled, the enigmatic web of chemical interactions unfolds before us, revealing the profound interconnectedness of nature</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Introduction Continued:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Chemistry weaves a tapestry of knowledge that connects the microcosm to the macrocosm, allowing us to understand the vastness of the universe down to the minute intricacies of life</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Its principles permeate various disciplines, shaping the modern world in countless ways</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> From uncovering the mysteries of DNA to the development of life-saving drugs, chemistry's impact is undeniable</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:lastRenderedPageBreak />
        <w:br />
        <w:t>In modern times, technological advancements driven by chemistry continue to revolutionize industries and shape our societies</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> From the realm of energy production to the development of cutting-edge materials, chemistry's contributions are boundless</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Understanding chemistry empowers us to grapple with global challenges such as climate change, food security, and resource depletion, seeking sustainable solutions that ensure a prosperous future for generations to come</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>At the forefront of scientific discovery, chemistry stands as a cornerstone of human knowledge, constantly evolving and adapting to new frontiers</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> As budding scientists and innovators, we're invited to embark on an enthralling expedition into the depths of chemistry, unravelling the intricate dance of atoms and molecules, and shaping the world around us for the better</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(5).Range.InsertXML($introXml) | Out-Null

# -----------------------------------------------------------------------
# 5) Summary body paragraph - full replace (run layout also changes).
# -----------------------------------------------------------------------
$summaryXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00653B50" w:rsidRDefault="005C06F8">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>Embark on a captivating journey through the enchanting realm of chemistry, exploring the union of atoms and molecules that orchestrate countless transformations in the world around us</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t xml:space="preserve"> Chemistry offers a comprehensive approach to comprehending the intricacies of chemical reactions, fascinating periodic trends, and the magic behind everyday occurrences</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t xml:space="preserve"> From the vibrant hues of chemical reactions to the intricate tapestry of interconnectedness, chemistry unveils a symphony of beauty and understanding that fosters a deeper appreciation for the wonders of science</w:t>
      </w:r>
      <w:r w:rsidR="00CA22DC">
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$d.Paragraphs(7).Range.InsertXML($summaryXml) | Out-Null

# -----------------------------------------------------------------------
# 6) Append a new trailing empty paragraph after the summary.
# -----------------------------------------------------------------------
$endRange = $d.Range($d.Content.End, $d.Content.End)
$emptyParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" />
'@
$endRange.InsertXML($emptyParaXml) | Out-Null

Write-Output "done"
